# Apply the "Add change log and roadmap" edit to tab_Tipo_tribunale.xlsx
#
# Summary of the change in the "SQL Results" sheet:
#   - Drop columns C:F (DESCRIZIONE/ORDINAMENTO/DATAINIZIOVALIDITA/DATAFINEVALIDITA
#     text columns + the always-empty D/E/F values), keeping only ID + DESCRIZIONE.
#   - Column A keeps the numeric ID values (1,2,3,4,9) instead of the old
#     row-number column; column B now holds the description text that used
#     to live in column C.
#   - "Avvocato/Notaio" becomes "Avvocato/Notaio/Uff.St.Civile".
#   - Column B autosizes/widens to fit the longer text, column A narrows.
#   - A thin border is added around the used range.

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("SQL Results")
$wsStatement = $wb.Worksheets.Item("SQL Statement")

# --- Rebuild "SQL Results" ---------------------------------------------

# Clear out the old C:F columns entirely (they disappear from the sheet).
$wsResults.Columns("C:F").Delete()

# Header row.
$wsResults.Range("A1").Value = "ID"
$wsResults.Range("B1").Value = "DESCRIZIONE"

# Data rows: id, description.
$data = @(
    @(1, "Tribunale civile Italiano"),
    @(2, "Sacra Rota"),
    @(3, "Tribunale Estero"),
    @(4, "Avvocato/Notaio/Uff.St.Civile"),
    @(9, "Non conosciuto")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $wsResults.Cells.Item($row, 1).Value = $data[$i][0]
    $wsResults.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column sizing: narrow ID column, widen (and best-fit) the description column.
$wsResults.Columns("A").ColumnWidth = 5
$wsResults.Columns("B").ColumnWidth = 23.85546875

# Bold header font + thin border around the whole used range (A1:B6).
$usedRange = $wsResults.Range("A1:B6")
$usedRange.Borders.LineStyle = 1
$usedRange.Borders.Weight = 2

$wsResults.Range("A1:B1").Font.Bold = $true

$wsResults.Range("C6").Select()

# --- "SQL Statement" sheet keeps the same text, just style bookkeeping --
$wsStatement.Range("A1").Value = "select * from DEC_TIPO_TRIBUNALE t"
